# Move module change password to Parent Tabs: change password requires log in
#
# This updates the localization workbook (tools/localize/localize.xlsx):
#  - Replace the placeholder English "log out confirm" string with an improved one
#  - Add English translations for 5 previously-untranslated description rows
#  - Shrink the row height of the 3 rows that no longer need 45pt
#  - Move the viewport / selection back up near the top of the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the English translation of the log-out confirmation text ---
$ws.Range("C99").Value = "Are You Sure To Log Out ?"

# --- Column C (English) translations for the "_des" rows (27-31) ---
$ws.Range("C28").Value = "Day off announcement to students"
$ws.Range("C27").Value = "Track the current location of students, bus comings"
$ws.Range("C29").Value = "Register the bus service, select the place, year, list of gủadians"
$ws.Range("C30").Value = "Change the register of bus service"
$ws.Range("C31").Value = "Register of guardians, who go with student on bus"

# Rows 27-29 no longer need the taller 45pt row (text is shorter now); row 31
# already was 30pt and stays that way.
$ws.Range("A27:C29").RowHeight = 30

# --- Restore the view near the top of the sheet / select C31 ---
$ws.Range("C31").Select()
